$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update stats for 2025-07 (row 20)
$ws.Range("B20").Value = 6172
$ws.Range("D20").Value = 5575345
$ws.Range("E20").Value = 903.3287427090085
$ws.Range("F20").Value = 6.615995854206247
$ws.Range("H20").Value = 26.1198238367778
